$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates ---
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Column F width (revert bestFit width back down to match columns C/D/E/G/H) ---
$ws.Columns.Item(6).ColumnWidth = 6.71

# --- Crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 379
$ws.Range("J14").Value = 432
$ws.Range("K14").Value = -12.268518518518
$ws.Range("L14").Value = -11.448598130841
$ws.Range("M14").Value = -21.369294605809
$ws.Range("N14").Value = -78.168202764977

# Row 15
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = 31
$ws.Range("E15").Value = -6.451612903225
$ws.Range("F15").Value = 118
$ws.Range("G15").Value = 117
$ws.Range("H15").Value = 0.8547008547
$ws.Range("I15").Value = 1471
$ws.Range("J15").Value = 1336
$ws.Range("K15").Value = 10.104790419161
$ws.Range("L15").Value = 12.633996937212
$ws.Range("M15").Value = 19.788273615635
$ws.Range("N15").Value = -49.846573474258

# Row 16
$ws.Range("C16").Value = 325
$ws.Range("D16").Value = 344
$ws.Range("E16").Value = -5.523255813953
$ws.Range("F16").Value = 1452
$ws.Range("G16").Value = 1398
$ws.Range("H16").Value = 3.862660944206
$ws.Range("I16").Value = 15639
$ws.Range("J16").Value = 12033
$ws.Range("K16").Value = 29.967589129892
$ws.Range("L16").Value = 34.935289042277
$ws.Range("M16").Value = -9.543640465035
$ws.Range("N16").Value = -79.480686470032

# Row 17
$ws.Range("C17").Value = 405
$ws.Range("D17").Value = 495
$ws.Range("E17").Value = -18.181818181818
$ws.Range("F17").Value = 1945
$ws.Range("G17").Value = 1910
$ws.Range("H17").Value = 1.832460732984
$ws.Range("I17").Value = 23270
$ws.Range("J17").Value = 20563
$ws.Range("K17").Value = 13.164421533822
$ws.Range("L17").Value = 24.899361279587
$ws.Range("M17").Value = 51.863212164719
$ws.Range("N17").Value = -37.660737248178

# Row 18
$ws.Range("C18").Value = 280
$ws.Range("D18").Value = 308
$ws.Range("E18").Value = -9.090909090909
$ws.Range("F18").Value = 1175
$ws.Range("G18").Value = 1179
$ws.Range("H18").Value = -0.339270568278
$ws.Range("I18").Value = 14000
$ws.Range("J18").Value = 11035
$ws.Range("K18").Value = 26.86905301314
$ws.Range("L18").Value = 1.699840185965
$ws.Range("M18").Value = -15.96134221742
$ws.Range("N18").Value = -84.363307383869

# Row 19
$ws.Range("C19").Value = 963
$ws.Range("D19").Value = 999
$ws.Range("E19").Value = -3.603603603603
$ws.Range("F19").Value = 4096
$ws.Range("G19").Value = 3916
$ws.Range("H19").Value = 4.596527068437
$ws.Range("I19").Value = 45960
$ws.Range("J19").Value = 33918
$ws.Range("K19").Value = 35.50327259862
$ws.Range("L19").Value = 44.079751716354
$ws.Range("M19").Value = 36.554060076655
$ws.Range("N19").Value = -39.787763657801

# Row 20
$ws.Range("C20").Value = 241
$ws.Range("D20").Value = 259
$ws.Range("E20").Value = -6.949806949806
$ws.Range("F20").Value = 1065
$ws.Range("G20").Value = 940
$ws.Range("H20").Value = 13.297872340425
$ws.Range("I20").Value = 12036
$ws.Range("J20").Value = 9077
$ws.Range("K20").Value = 32.598876280709
$ws.Range("L20").Value = 50.676014021031
$ws.Range("M20").Value = 30.287941112795
$ws.Range("N20").Value = -87.937220629798

# Row 21
$ws.Range("C21").Value = 2248
$ws.Range("D21").Value = 2441
$ws.Range("E21").Value = -7.906595657517
$ws.Range("F21").Value = 9876
$ws.Range("G21").Value = 9485
$ws.Range("H21").Value = 4.12229836584
$ws.Range("I21").Value = 112755
$ws.Range("J21").Value = 88394
$ws.Range("K21").Value = 27.559562866257
$ws.Range("L21").Value = 31.710821418559
$ws.Range("M21").Value = 20.110571392049
$ws.Range("N21").Value = -70.625550339452

# Row 22
$ws.Range("C22").Value = 36
$ws.Range("D22").Value = 59
$ws.Range("E22").Value = -38.983050847457
$ws.Range("F22").Value = 193
$ws.Range("G22").Value = 205
$ws.Range("H22").Value = -5.853658536585
$ws.Range("I22").Value = 2056
$ws.Range("J22").Value = 1521
$ws.Range("K22").Value = 35.174227481919
$ws.Range("L22").Value = 28.259513412351
$ws.Range("M22").Value = 8.496042216358

# Row 23
$ws.Range("C23").Value = 122
$ws.Range("D23").Value = 131
$ws.Range("E23").Value = -6.870229007633
$ws.Range("F23").Value = 460
$ws.Range("G23").Value = 467
$ws.Range("H23").Value = -1.498929336188
$ws.Range("I23").Value = 5305
$ws.Range("J23").Value = 4906
$ws.Range("K23").Value = 8.132898491642
$ws.Range("L23").Value = 15.956284153005
$ws.Range("M23").Value = 40.977943130481

# Row 24
$ws.Range("C24").Value = 2324
$ws.Range("D24").Value = 2014
$ws.Range("E24").Value = 15.392254220456
$ws.Range("F24").Value = 8973
$ws.Range("G24").Value = 7913
$ws.Range("H24").Value = 13.39567799823
$ws.Range("I24").Value = 102820
$ws.Range("J24").Value = 75285
$ws.Range("K24").Value = 36.574350800292
$ws.Range("L24").Value = 42.100971571513
$ws.Range("M24").Value = 40.924603555324

# Row 25
$ws.Range("C25").Value = 687
$ws.Range("D25").Value = 792
$ws.Range("E25").Value = -13.257575757575
$ws.Range("F25").Value = 3214
$ws.Range("G25").Value = 3105
$ws.Range("H25").Value = 3.510466988727
$ws.Range("I25").Value = 36967
$ws.Range("J25").Value = 32178
$ws.Range("K25").Value = 14.882839206911
$ws.Range("L25").Value = 24.021203073103
$ws.Range("M25").Value = -10.504527187333

# Row 26
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 51
$ws.Range("E26").Value = -13.725490196078
$ws.Range("F26").Value = 180
$ws.Range("G26").Value = 207
$ws.Range("H26").Value = -13.043478260869
$ws.Range("I26").Value = 2347
$ws.Range("J26").Value = 2194
$ws.Range("K26").Value = 6.97356426618
$ws.Range("L26").Value = 14.767726161369

# Row 27
$ws.Range("C27").Value = 107
$ws.Range("D27").Value = 102
$ws.Range("E27").Value = 4.901960784313
$ws.Range("F27").Value = 444
$ws.Range("G27").Value = 412
$ws.Range("H27").Value = 7.766990291262
$ws.Range("I27").Value = 4671
$ws.Range("J27").Value = 4385
$ws.Range("K27").Value = 6.522234891676
$ws.Range("L27").Value = 35.509138381201

# Row 28
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 37
$ws.Range("E28").Value = -56.756756756756
$ws.Range("F28").Value = 96
$ws.Range("G28").Value = 135
$ws.Range("H28").Value = -28.888888888888
$ws.Range("I28").Value = 1438
$ws.Range("J28").Value = 1688
$ws.Range("K28").Value = -14.810426540284
$ws.Range("L28").Value = -14.709371293001
$ws.Range("M28").Value = -10.293200249532
$ws.Range("N28").Value = -73.166635566337

# Row 29
$ws.Range("C29").Value = 14
$ws.Range("D29").Value = 25
$ws.Range("E29").Value = -44
$ws.Range("F29").Value = 79
$ws.Range("G29").Value = 109
$ws.Range("H29").Value = -27.522935779816
$ws.Range("I29").Value = 1183
$ws.Range("J29").Value = 1403
$ws.Range("K29").Value = -15.680684248039
$ws.Range("L29").Value = -13.901018922853
$ws.Range("M29").Value = -10.514372163388
$ws.Range("N29").Value = -75.420735507999

# Row 30
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = -33.333333333333
$ws.Range("F30").Value = 39
$ws.Range("G30").Value = 39
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 556
$ws.Range("J30").Value = 481
$ws.Range("K30").Value = 15.592515592515
$ws.Range("L30").Value = 128.80658436214
